# Update the May 2023 DTR report table: rendered hours dropped for several
# members, which re-ranks the (descending by hours) table and flips some
# rows from Complete to Incomplete.
#
# Rather than moving w:tr nodes around, we simply rewrite every data cell
# with the new target content, row by row, which yields an OOXML result
# equivalent to the reordered table described by the diff.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (header) is unchanged: Name | Position | Total Rendered Hours | Remarks

$rows = @(
    @("BINONDO, KYLA D.", "Editorial Board Managing Director", "199:33:07", "Complete"),
    @("SISMAR, KARREN MARIE B.", "Senior Staff Photojournalist", "149:36:27", "Complete"),
    @("DELOS REYES, HASNA ALTHEA M.", "Editorial Board Editor in Chief", "130:52:33", "Complete"),
    @("SABANAL, JUVYL T.", "Senior Staff Layout Artist", "100:02:42", "Complete"),
    @("BARRIENTOS, JOHN CLEISTER C.", "Senior Staff Photojournalist", "93:41:04", "Complete"),
    @("BINONDO, KYZEN D.", "Editorial Board Feature Editor", "82:46:44", "Complete"),
    @("DEIMOS, CHRISTIAN JACOB B.", "Senior Staff Photojournalist", "75:24:01", "Complete"),
    @("EYAO, ADRIENNE C.", "Senior Staff Layout Artist", "71:35:59", "Complete"),
    @("AMPO-ON, SARC FRANCIS ADRIANNE  T.", "Senior Staff Cartoonist", "62:58:19", "Complete"),
    @("HONTIVEROS, MARK DENVER  Y.", "Senior Staff Cartoonist", "56:25:31", "Complete"),
    @("RACAZA, DAVE N.", "Senior Staff Writer", "51:11:27", "Complete"),
    @("NUÑEZ, LEN D.", "Senior Staff Layout Artist", "51:10:32", "Complete"),
    @("ILLUT, NICHOLS JOHN M.", "Senior Staff Photojournalist", "38:57:45", "Incomplete"),
    @("LANZADERAS, MARIE CHASTINE V.", "Senior Staff Layout Artist", "32:45:48", "Incomplete"),
    @("BUGHAO, NINETTE ANN C.", "Editorial Board Art Editor", "26:31:50", "Incomplete"),
    @("SOLON, ADAM L.", "Editorial Board Photo Editor", "25:46:45", "Incomplete"),
    @("LAURONILLA, ARGYLE JOSEPH  M.", "Senior Staff Writer", "18:12:44", "Incomplete"),
    @("GIMENEZ, CATHYRENE A.", "Editorial Board Finance Manager", "16:41:35", "Incomplete"),
    @("TORRES, ANGELA MAE S.", "Editorial Board Planning and Research Director", "10:59:52", "Incomplete"),
    @("DELA CRUZ, REMUEL B.", "Senior Staff Cartoonist", "4:57:36", "Incomplete"),
    @("VELEZ, TRIXIA GLENN B.", "Senior Staff Writer", "1:44:08", "Incomplete")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowIndex = $i + 2   # row 1 is the header
    $vals = $rows[$i]
    for ($col = 1; $col -le 4; $col++) {
        $t.Cell($rowIndex, $col).Range.Text = $vals[$col - 1]
    }
}

Write-Output "Done updating table rows."
